$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: rotate "Gemma 3 Fine-tuned"(P1) -> N1, "Gemma 3 RAG"(O1) -> P1,
#     keep "Gemma 3 Embed" label at O1 (new position) ---
# Capture target formats first (before any value/format writes change things):
$ws.Range("P1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N1").Value = "Gemma 3 Fine-tuned"
$ws.Range("O1").Value = "Gemma 3 Embed"
$ws.Range("P1").Value = "Gemma 3 RAG"

# --- Data rows 2-8: N <- old P (value & style), P <- old O (value, style already matches),
#     O <- brand-new "Gemma 3 Embed" test values (style already matches, just set value) ---

$rows = 2,3,4,5,6,7,8
$oldP = @{}
$oldO = @{}
foreach ($r in $rows) {
    $oldP[$r] = $ws.Range("P$r").Value2
    $oldO[$r] = $ws.Range("O$r").Value2
}

# New "Gemma 3 Embed" values to populate column O
$newO = @{
    2 = 0.579
    3 = 0.6
    4 = 0.75
    5 = 0.7
    6 = 0.75
    7 = 0.55
    8 = 0.65
}

foreach ($r in $rows) {
    # N gets old P's format (copy BEFORE P's own format is touched)
    $ws.Range("P$r").Copy()
    $ws.Range("N$r").PasteSpecial(-4122)   # xlPasteFormats

    # P gets old O's format (O's per-cell style never changes, so this is a no-op in practice,
    # but keeps things robust/explicit)
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)   # xlPasteFormats

    # Now assign the actual values
    $ws.Range("N$r").Value = $oldP[$r]
    $ws.Range("P$r").Value = $oldO[$r]
    $ws.Range("O$r").Value = $newO[$r]
}

$excel.CutCopyMode = 0

# --- Column widths (best effort given COM ColumnWidth integer-pixel rounding) ---
$ws.Columns("N").ColumnWidth = 23.142857142857142
$ws.Columns("O").ColumnWidth = 18.428571428571427
$ws.Columns("P").ColumnWidth = 17.428571428571427

# --- Sheet view: scroll back to A1 (drop topLeftCell="C1") and move selection ---
$ws.Range("A1").Select()
$ws.Range("I16").Select()
